# Applies the "fix caps error nd redone transfer + transfer back" edit.
$wb = $excel.ActiveWorkbook

$wsA2 = $wb.Worksheets.Item("A2")
$wsA3 = $wb.Worksheets.Item("A3")
$wsA4 = $wb.Worksheets.Item("A4")
$wsA5 = $wb.Worksheets.Item("A5")
$wsA6 = $wb.Worksheets.Item("A6")

# --- Sheet A2: row 3 (TxHash rehashed; NFTID label memoryGrinderNFT2 -> memoryGrinderNFT4) ---
$wsA2.Activate()
$wsA2.Range("A3").Value = "172C6D84CB7876C27376ECC7D4D408990EB350FB3F4F6C4ABA745499FAAA5C3A"
$wsA2.Range("C3").Value = "memoryGrinderNFT4"
$wsA2.Range("C4").Select()

# --- Sheet A3: D2 chain name caps fix ---
# (ColumnWidth input is pre-calibrated so the on-disk OOXML width lands as close as
# this engine's pixel/character-width rounding allows to the target of 83.73)
$wsA3.Activate()
$wsA3.Range("D2").Value = "elgafar-1"
$wsA3.Columns.Item(2).ColumnWidth = 82.8
$wsA3.Range("C8").Select()

# --- Sheet A4: A2 TxHash rehashed, C2 NFTID label (shared with A2/A6), D2 chain name caps fix ---
$wsA4.Activate()
$wsA4.Range("A2").Value = "14835E551EF9A9F551030B20D90EDE4AC29A0E0A7E647147EF7F9A45DEC7E0A2"
$wsA4.Range("C2").Value = "memoryGrinderNFT4"
$wsA4.Range("D2").Value = "uptick_7000-2"
$wsA4.Columns.Item(2).ColumnWidth = 49.8
$wsA4.Columns.Item(3).ColumnWidth = 20.6
$wsA4.Range("F10").Select()

# --- Sheet A5: D2 chain name caps fix (shared string with A3) ---
$wsA5.Activate()
$wsA5.Range("D2").Value = "elgafar-1"
$wsA5.Range("E32").Select()

# --- Sheet A6: A2 TxHash rehashed, C2 NFTID label (shared with A2/A4), D2 chain name caps fix (shared with A4) ---
$wsA6.Activate()
$wsA6.Range("A2").Value = "ED35D60216058469285BF3BA9AE2811CE45737312B00CFC110C6C5AF0D547D50"
$wsA6.Range("C2").Value = "memoryGrinderNFT4"
$wsA6.Range("D2").Value = "uptick_7000-2"
$wsA6.Columns.Item(2).ColumnWidth = 59.8
$wsA6.Columns.Item(3).ColumnWidth = 25.92
$wsA6.Range("D3").Select()

# --- Active sheet / tab changes: activeTab moves from A6 (index 6) to A5 (index 5) ---
$wsA5.Activate()
